# Add a new "Hat" price entry ("Gjøgler pacman lilla" @ 5200) to the price
# list, keeping it in alphabetical order, then turn the A1:B28 range into a
# proper Excel Table (ListObject) so new prices can be added/filtered easily.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The list is sorted alphabetically; "Gjøgler pacman lilla" belongs right
# before "Gjøglertiss", which currently lives on row 11. Insert a fresh row
# there (pushing everything below it down by one) and fill in the new item.
$ws.Rows.Item(11).Insert()
$ws.Range("A11").Value = "Gjøgler pacman lilla"
$ws.Range("A11").Font.Name = "Arial"
$ws.Range("A11").Font.Size = 10
$ws.Range("B11").Value = 5200

# Turn the whole range (now A1:B28, header included) into a real table so it
# can be sorted/filtered and new rows appended easily.
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:B28"), $null, 1)
$tbl.Name = "Table1"

# Match the author's final selection.
[void]$ws.Range("F7").Select()
